# Contract Whist Scorecard / Colors.xlsx
# "Tightened up awards - especially sync and defaulted awards
#  Removed sync from high scores dashboard from game summary
#  Made sync OK for landscape"
#
# The bulk of the real content edit lives on the "Views" worksheet: a
# bunch of previously-blank "Specified/Landscape/Portrait" tracking cells
# get marked "Yes" (spec now defined for those items), and a handful of
# manual note cells in column J ("Check ... with Jack", "Shield
# definitions", ...) get cleared out now that those items are resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Views")

# Cells that move from blank to "Yes" (spec/landscape/portrait now done).
$yesCells = @(
    "D3", "E3",
    "C10",
    "D12", "E12",
    "B13", "C13",
    "D18", "E18",
    "D19", "E19",
    "D20", "E20",
    "D21", "E21",
    "D22", "E22",
    "D23", "E23",
    "D24", "E24",
    "D25", "E25",
    "C27", "D27", "E27"
)
foreach ($cellRef in $yesCells) {
    $ws.Range($cellRef).Value = "Yes"
}

# Manual "waiting on Jack" style notes in column J that are no longer
# needed now that the corresponding rows are resolved - clear back to
# blank (not just empty string) so they fall back to whatever the shared
# formula produces, or stay genuinely empty where there was no formula.
$clearCells = @("J11", "J14", "J15", "J18", "J23", "J25", "J29")
foreach ($cellRef in $clearCells) {
    $ws.Range($cellRef).ClearContents()
}

# Selection on the Views sheet moved from C10 to K29.
$ws.Range("K29").Select() | Out-Null

Write-Output "Views sheet updated"
